$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Date and Publisher values
$ws.Range("B8").Value = "2021-12-22T21:26:07+01:00"
$ws.Range("B9").Value = "Forschungsgruppe Digital Health"

# Insert two new rows after row 10 (Contact row), shifting everything below down
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# New row 11 duplicates the Contact property (matches source content)
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "No display for ContactDetail"

# New row 12 is the Jurisdiction property
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = "Germany"
